# Change in gender slide title
# Slide 7 ("Does gender affect sleep?") title text is updated to the
# new, more specific research question.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$titleShape = $s.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Does age & gender affects sleep efficiency & sleep duration ?"
